# "Changes in b suite" — flip the Results value for TestCase_E6 (row 6) on
# the "Test Cases" sheet from "FAIL" to "SKIP" (matches every other row in
# the Results column, which already reads "SKIP").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Range("E6").Value = "SKIP"
